$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2018 LEAVE CREDITS")
$ws.Activate()
$win = $excel.ActiveWindow
$win.SplitRow = 57
$win.SplitColumn = 0
$ws.Range("I90").Select()
$win.ScrollRow = 68
$win.ScrollColumn = 1
